# Apply crime data update for 2023-12-17
# Each worksheet below has one or more YTD cumulative cell values incremented
# to reflect a newly recorded incident as of 2023-12-17.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 116
$ws.Range("D3").Value = 139
$ws.Range("E3").Value = 148
$ws.Range("H3").Value = 158
$ws.Range("J3").Value = 239
$ws.Range("H4").Value = 14
$ws.Range("B6").Value = 384
$ws.Range("C6").Value = 490
$ws.Range("D6").Value = 424
$ws.Range("E6").Value = 489
$ws.Range("H6").Value = 448
$ws.Range("I6").Value = 508
$ws.Range("J6").Value = 427
$ws.Range("B7").Value = 518
$ws.Range("C7").Value = 647
$ws.Range("D7").Value = 664
$ws.Range("E7").Value = 723
$ws.Range("H7").Value = 732
$ws.Range("I7").Value = 847
$ws.Range("J7").Value = 812

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("D3").Value = 9
$ws.Range("E6").Value = 54
$ws.Range("D7").Value = 49
$ws.Range("E7").Value = 67

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I6").Value = 21
$ws.Range("I7").Value = 37

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("H4").Value = 2
$ws.Range("B6").Value = 32
$ws.Range("B7").Value = 37
$ws.Range("H7").Value = 47

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("E8").Value = 54
$ws.Range("I8").Value = 42
$ws.Range("B28").Value = 37
$ws.Range("H28").Value = 47
$ws.Range("D32").Value = 49
$ws.Range("E32").Value = 67
$ws.Range("I36").Value = 37
$ws.Range("H47").Value = 25
$ws.Range("J47").Value = 17
$ws.Range("C51").Value = 3
$ws.Range("D53").Value = 75
$ws.Range("E53").Value = 88
$ws.Range("H53").Value = 103
$ws.Range("I53").Value = 126
$ws.Range("J53").Value = 125
$ws.Range("H70").Value = 17
$ws.Range("J76").Value = 17
$ws.Range("D85").Value = 6
$ws.Range("B98").Value = 518
$ws.Range("C98").Value = 647
$ws.Range("D98").Value = 664
$ws.Range("E98").Value = 723
$ws.Range("H98").Value = 732
$ws.Range("I98").Value = 847
$ws.Range("J98").Value = 812

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("I2").Value = 14
$ws.Range("E3").Value = 17
$ws.Range("H3").Value = 23
$ws.Range("J3").Value = 38
$ws.Range("D6").Value = 46
$ws.Range("J6").Value = 62
$ws.Range("D7").Value = 75
$ws.Range("E7").Value = 88
$ws.Range("H7").Value = 103
$ws.Range("I7").Value = 126
$ws.Range("J7").Value = 125

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("J3").Value = 3
$ws.Range("J7").Value = 17

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("C4").Value = 2
$ws.Range("C5").Value = 3
$ws.Range("D4").Value = 4

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("D5").Value = 6

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("H5").Value = 15
$ws.Range("J5").Value = 4
$ws.Range("H6").Value = 25
$ws.Range("J6").Value = 17

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("H4").Value = 12
$ws.Range("H5").Value = 17

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("E5").Value = 42
$ws.Range("I5").Value = 32
$ws.Range("E6").Value = 54
$ws.Range("I6").Value = 42
